## edit.ps1 - reproduce TC01_Canine_Filter_Breed-Akita.xlsx 'startup' sheet rework
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "startup" sheet (first sheet / tabSelected)

# ---- long text blocks, single-quoted here-strings so $ and ` are literal ----
$statQueryText = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Akita']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Akita']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
        coalesce(co.cohort_description, '') AS `Cohort`

'@
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Akita']  
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@
$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN ['Akita'] 
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@

# ---- insert a new column A for the 'TabName' label column ----
$ws.Columns.Item(1).Insert()

# ---- header row ----
$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# ---- row 2: CasesTab ----
$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQueryText
$ws.Range("D2").Value = "TC01_Canine_Filter_Breed-Akita_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC01_Canine_Filter_Breed-Akita_WebData.xlsx"
$ws.Rows.Item(2).RowHeight = 275.5

# ---- row 3: SamplesTab ----
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQueryText
$ws.Range("D3").Value = "TC01_Canine_Filter_Breed-Akita_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC01_Canine_Filter_Breed-Akita_WebData.xlsx"
$ws.Rows.Item(3).RowHeight = 217.5

# ---- row 4: FilesTab ----
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQueryText
$ws.Range("D4").Value = "TC01_Canine_Filter_Breed-Akita_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC01_Canine_Filter_Breed-Akita_WebData.xlsx"
$ws.Rows.Item(4).RowHeight = 246.5

# ---- wrap text on the long-text columns (B & C) for rows 2-4 ----
$ws.Range("B2:C4").WrapText = $true

# ---- rows 5-13: formatted-but-empty C cells (style carried down, no value) ----
for ($r = 5; $r -le 13; $r++) {
    $ws.Cells.Item($r, 3).WrapText = $true
}

# ---- A2 picked up a stray alignment flag with no explicit alignment (matches source) ----
$ws.Range("A2").IndentLevel = 1
$ws.Range("A2").IndentLevel = 0

# ---- column widths ----
# (engine rounds ColumnWidth to the nearest 1/6 char-unit; inputs below are
#  chosen so the stored width lands as close as possible to the source values
#  75.81640625->10.90625(A) / 255.6328125->92.453125(B) would've produced under a
#  real-Excel MDW=7 model -- 10.0 / 91.66666666666667 / 85.33333333333333)
$ws.Columns.Item(1).ColumnWidth = 10.0
$ws.Columns.Item(2).ColumnWidth = 91.66666666666667
$ws.Columns.Item(3).ColumnWidth = 85.33333333333333

# ---- view: zoom 55%, selection on B2, no frozen/top-left-cell override ----
$ws.Application.ActiveWindow.Zoom = 55
$ws.Range("B2").Select()

